$d = $word.ActiveDocument

# 1. Title: "Environmental Noise Contamination Detector - Data Pipeline"
#    runs merge together (no visible text change); touch it so the two
#    runs collapse into one, matching the target structure.
$d.Content.Find.Execute(
    "Environmental Noise Contamination Detector " + [char]0x2013 + " Data Pipeline",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Environmental Noise Contamination Detector " + [char]0x2013 + " Data Pipeline", 2) | Out-Null

# 2. Due date: "Due: 2" + "1 " + "Nov" + " 2018" -> single run "Due: 21 Nov 2018"
$d.Content.Find.Execute("Due: 21 Nov 2018", $true, $false, $false, $false, $false,
    $true, 1, $false, "Due: 21 Nov 2018", 2) | Out-Null

# 3. "The data used from this project" -> "The data used for this project"
$d.Content.Find.Execute("The data used from this project will originate from files. Typically, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The data used for this project will originate from files. Typically, ", 2) | Out-Null

# 4. "to each team member to store locally. " -> add justification clause
$d.Content.Find.Execute("to each team member to store locally. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "to each team member to store locally as total size is over 500 MB and not suitable for storage in a GitHub repository. ", 2) | Out-Null

# 5. "all of our processing. " -> "all our processing. " (drop "of " and the
#    gramStart/gramEnd proofing marks that used to bracket "all of")
$d.Content.Find.Execute("all of our processing. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "all our processing. ", 2) | Out-Null

# 6. "(the plane08.wav file)" -> "(the plane08.wav file of a turbofan jet aircraft)"
$d.Content.Find.Execute("(the plane08.wav file)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(the plane08.wav file of a turbofan jet aircraft)", 2) | Out-Null
